# Updated symbol list on Sun Dec 25 10:42:17 UTC 2022 with GitHub Actions
#
# Applies the latest scrape refresh to the crypto price sheet: a batch of
# small price (column D) corrections, plus three rows (41-43) whose coins
# got re-ranked (BKEXToken / CEJI / KickToken shuffled), and a couple of
# "Volume(1h)" label tweaks (column E) that follow from the re-rank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Numeric-looking values live in text/inline-string cells in this sheet
    # (e.g. "0.05980" with a significant trailing zero). Force the cell to
    # Text format first so Excel doesn't silently convert it to a Number
    # and strip the trailing zero / change precision, then restore the
    # cell's original "Normal" style so only the value changes (no leftover
    # explicit number-format override stays on the cell).
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value2 = $value
    $cell.Style = "Normal"
}

# ---- simple price (column D) corrections ----
Set-TextValue "D2"  "244.71"
Set-TextValue "D3"  "23.06"
Set-TextValue "D5"  "0.05980"
Set-TextValue "D6"  "3.391"
Set-TextValue "D8"  "0.9276"
Set-TextValue "D9"  "0.1432"
Set-TextValue "D10" "0.07435"
Set-TextValue "D11" "0.03370"
Set-TextValue "D13" "0.09340"
Set-TextValue "D14" "3.935"
Set-TextValue "D15" "0.001591"
Set-TextValue "D17" "0.0005943"
Set-TextValue "D18" "0.005638"
Set-TextValue "D19" "0.004155"
Set-TextValue "D20" "0.0009837"
Set-TextValue "D23" "6.458"
Set-TextValue "D40" "0.03942"

# ---- rows 41-43: BKEXToken / CEJI / KickToken re-ranked ----
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1075"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002711"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003021"
$ws.Range("E43").Value = "42KickTokenKICK"

# ---- row 44: LocalTraders price + volume-label refresh ----
Set-TextValue "D44" "0.007487"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"

# ---- remaining simple price (column D) corrections ----
Set-TextValue "D45" "0.00005132"
Set-TextValue "D47" "0.0005803"
Set-TextValue "D49" "0.002255"
